$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 (bold, bordered, centered) onto the new
# I1 / J1 header cells so they match the existing header row style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New "I0" / "IF" data columns, keyed by row number -> @(I-value, J-value)
$ijData = @{
    2 = @(8, 8)
    3 = @(6, 7)
    4 = @(6, 6)
    5 = @(8, 8)
    6 = @(8, 8)
    7 = @(7, 7)
    8 = @(7, 7)
    9 = @(8, 8)
    10 = @(6, 6)
    11 = @(7, 7)
    12 = @(7, 7)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(7, 8)
    16 = @(7, 7)
    17 = @(8, 8)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(6, 6)
    21 = @(9, 9)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(6, 6)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(8, 8)
    32 = @(8, 8)
    33 = @(8, 8)
    34 = @(9, 9)
    35 = @(8, 8)
    36 = @(10, 10)
    37 = @(7, 8)
    38 = @(8, 8)
    39 = @(8, 8)
    40 = @(7, 8)
    41 = @(7, 7)
    42 = @(7, 7)
    43 = @(6, 6)
    44 = @(10, 10)
    45 = @(7, 7)
    46 = @(8, 8)
    47 = @(7, 7)
    48 = @(6, 6)
    49 = @(9, 9)
    50 = @(6, 7)
    51 = @(8, 8)
    52 = @(6, 7)
    53 = @(10, 10)
    54 = @(8, 8)
    55 = @(8, 8)
    56 = @(6, 6)
    57 = @(4, 5)
    58 = @(8, 8)
    59 = @(6, 6)
    60 = @(7, 7)
    61 = @(6, 6)
    62 = @(6, 6)
    63 = @(7, 7)
    64 = @(6, 6)
    65 = @(6, 6)
    66 = @(9, 9)
    67 = @(8, 9)
    68 = @(7, 7)
    69 = @(6, 6)
    70 = @(5, 5)
    71 = @(4, 5)
    72 = @(6, 6)
}

foreach ($r in $ijData.Keys) {
    $pair = $ijData[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
